$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------
# Sheet "Schedule": rewrite row 2 and append rows 3-5
# (run 146 re-optimised the pumping schedule over a longer horizon)
# ----------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Cells.Item(2, "A").Value = 46043
$schedule.Cells.Item(2, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(2, "B").Value = 46043.16666666666
$schedule.Cells.Item(2, "B").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(2, "C").Value = 4
$schedule.Cells.Item(2, "D").Value = 15.12
$schedule.Cells.Item(2, "E").Value = 506.74358475
$schedule.Cells.Item(2, "F").Value = 33.51478735119048

$schedule.Cells.Item(3, "A").Value = 46043.29166666666
$schedule.Cells.Item(3, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(3, "B").Value = 46043.66666666666
$schedule.Cells.Item(3, "B").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(3, "C").Value = 9
$schedule.Cells.Item(3, "D").Value = 34.02
$schedule.Cells.Item(3, "E").Value = -160.897698
$schedule.Cells.Item(3, "F").Value = -4.729503174603174

$schedule.Cells.Item(4, "A").Value = 46043.875
$schedule.Cells.Item(4, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(4, "B").Value = 46044.10416666666
$schedule.Cells.Item(4, "B").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(4, "C").Value = 5.5
$schedule.Cells.Item(4, "D").Value = 20.79
$schedule.Cells.Item(4, "E").Value = 635.03264175
$schedule.Cells.Item(4, "F").Value = 30.54510061327561

$schedule.Cells.Item(5, "A").Value = 46044.27083333334
$schedule.Cells.Item(5, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(5, "B").Value = 46044.66666666666
$schedule.Cells.Item(5, "B").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(5, "C").Value = 9.5
$schedule.Cells.Item(5, "D").Value = 35.91
$schedule.Cells.Item(5, "E").Value = -21.09909750000001
$schedule.Cells.Item(5, "F").Value = -0.5875549289891399

# ----------------------------------------------------------------
# Sheet "Detailed": update Pump_Status / Price / Type cells
# ----------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Cells.Item(10, "E").Value = "OFF"

$detailed.Cells.Item(11, "E").Value = "OFF"

$detailed.Cells.Item(12, "E").Value = "OFF"

$detailed.Cells.Item(13, "E").Value = "OFF"

$detailed.Cells.Item(14, "E").Value = "OFF"

$detailed.Cells.Item(15, "E").Value = "OFF"

$detailed.Cells.Item(25, "B").Value = -15.89865

$detailed.Cells.Item(26, "B").Value = -22.10072

$detailed.Cells.Item(27, "B").Value = -23.5
$detailed.Cells.Item(27, "C").Value = "historical"

$detailed.Cells.Item(28, "B").Value = -24.41017
$detailed.Cells.Item(28, "C").Value = "historical"

$detailed.Cells.Item(29, "B").Value = -23.5
$detailed.Cells.Item(29, "C").Value = "historical"

$detailed.Cells.Item(30, "B").Value = -12.01
$detailed.Cells.Item(30, "C").Value = "historical"

$detailed.Cells.Item(31, "B").Value = -11.01
$detailed.Cells.Item(31, "C").Value = "historical"

$detailed.Cells.Item(32, "B").Value = -8.52

$detailed.Cells.Item(33, "B").Value = -5.85572

$detailed.Cells.Item(34, "B").Value = -5.45297

$detailed.Cells.Item(35, "B").Value = 0.00005

$detailed.Cells.Item(36, "B").Value = 36.25

$detailed.Cells.Item(37, "B").Value = 48.74546

$detailed.Cells.Item(39, "B").Value = 65

$detailed.Cells.Item(41, "B").Value = 73.45771000000001

$detailed.Cells.Item(42, "B").Value = 79.95

$detailed.Cells.Item(43, "B").Value = 72.57914

$detailed.Cells.Item(44, "B").Value = 64.8901
$detailed.Cells.Item(44, "E").Value = "ON"

$detailed.Cells.Item(45, "B").Value = 59.34928
$detailed.Cells.Item(45, "E").Value = "ON"

$detailed.Cells.Item(46, "E").Value = "ON"

$detailed.Cells.Item(47, "B").Value = 59.01655
$detailed.Cells.Item(47, "E").Value = "ON"

$detailed.Cells.Item(48, "B").Value = 57.31
$detailed.Cells.Item(48, "E").Value = "ON"

$detailed.Cells.Item(49, "B").Value = 62.04883
$detailed.Cells.Item(49, "E").Value = "ON"

# ----------------------------------------------------------------
# Sheet "Detailed": append new rows 50-97 (second day of the run)
# ----------------------------------------------------------------
$detailed.Cells.Item(50, "A").Value = 46044
$detailed.Cells.Item(50, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(50, "B").Value = 57.31
$detailed.Cells.Item(50, "C").Value = "forecast"
$detailed.Cells.Item(50, "D").Value = 46044
$detailed.Cells.Item(50, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(50, "E").Value = "ON"

$detailed.Cells.Item(51, "A").Value = 46044.02083333334
$detailed.Cells.Item(51, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(51, "B").Value = 62.40077
$detailed.Cells.Item(51, "C").Value = "forecast"
$detailed.Cells.Item(51, "D").Value = 46044
$detailed.Cells.Item(51, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(51, "E").Value = "ON"

$detailed.Cells.Item(52, "A").Value = 46044.04166666666
$detailed.Cells.Item(52, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(52, "B").Value = 57.31
$detailed.Cells.Item(52, "C").Value = "forecast"
$detailed.Cells.Item(52, "D").Value = 46044
$detailed.Cells.Item(52, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(52, "E").Value = "ON"

$detailed.Cells.Item(53, "A").Value = 46044.0625
$detailed.Cells.Item(53, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(53, "B").Value = 57.31
$detailed.Cells.Item(53, "C").Value = "forecast"
$detailed.Cells.Item(53, "D").Value = 46044
$detailed.Cells.Item(53, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(53, "E").Value = "ON"

$detailed.Cells.Item(54, "A").Value = 46044.08333333334
$detailed.Cells.Item(54, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(54, "B").Value = 57.06
$detailed.Cells.Item(54, "C").Value = "forecast"
$detailed.Cells.Item(54, "D").Value = 46044
$detailed.Cells.Item(54, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(54, "E").Value = "ON"

$detailed.Cells.Item(55, "A").Value = 46044.10416666666
$detailed.Cells.Item(55, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(55, "B").Value = 57.31
$detailed.Cells.Item(55, "C").Value = "forecast"
$detailed.Cells.Item(55, "D").Value = 46044
$detailed.Cells.Item(55, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(55, "E").Value = "OFF"

$detailed.Cells.Item(56, "A").Value = 46044.125
$detailed.Cells.Item(56, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(56, "B").Value = 64.89
$detailed.Cells.Item(56, "C").Value = "forecast"
$detailed.Cells.Item(56, "D").Value = 46044
$detailed.Cells.Item(56, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(56, "E").Value = "OFF"

$detailed.Cells.Item(57, "A").Value = 46044.14583333334
$detailed.Cells.Item(57, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(57, "B").Value = 64.89
$detailed.Cells.Item(57, "C").Value = "forecast"
$detailed.Cells.Item(57, "D").Value = 46044
$detailed.Cells.Item(57, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(57, "E").Value = "OFF"

$detailed.Cells.Item(58, "A").Value = 46044.16666666666
$detailed.Cells.Item(58, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(58, "B").Value = 64.89
$detailed.Cells.Item(58, "C").Value = "forecast"
$detailed.Cells.Item(58, "D").Value = 46044
$detailed.Cells.Item(58, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(58, "E").Value = "OFF"

$detailed.Cells.Item(59, "A").Value = 46044.1875
$detailed.Cells.Item(59, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(59, "B").Value = 66.268
$detailed.Cells.Item(59, "C").Value = "forecast"
$detailed.Cells.Item(59, "D").Value = 46044
$detailed.Cells.Item(59, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(59, "E").Value = "OFF"

$detailed.Cells.Item(60, "A").Value = 46044.20833333334
$detailed.Cells.Item(60, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(60, "B").Value = 66.87584
$detailed.Cells.Item(60, "C").Value = "forecast"
$detailed.Cells.Item(60, "D").Value = 46044
$detailed.Cells.Item(60, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(60, "E").Value = "OFF"

$detailed.Cells.Item(61, "A").Value = 46044.22916666666
$detailed.Cells.Item(61, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(61, "B").Value = 77.74561
$detailed.Cells.Item(61, "C").Value = "forecast"
$detailed.Cells.Item(61, "D").Value = 46044
$detailed.Cells.Item(61, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(61, "E").Value = "OFF"

$detailed.Cells.Item(62, "A").Value = 46044.25
$detailed.Cells.Item(62, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(62, "B").Value = 74.33651
$detailed.Cells.Item(62, "C").Value = "forecast"
$detailed.Cells.Item(62, "D").Value = 46044
$detailed.Cells.Item(62, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(62, "E").Value = "OFF"

$detailed.Cells.Item(63, "A").Value = 46044.27083333334
$detailed.Cells.Item(63, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(63, "B").Value = 57.06
$detailed.Cells.Item(63, "C").Value = "forecast"
$detailed.Cells.Item(63, "D").Value = 46044
$detailed.Cells.Item(63, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(63, "E").Value = "ON"

$detailed.Cells.Item(64, "A").Value = 46044.29166666666
$detailed.Cells.Item(64, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(64, "B").Value = 26.46921
$detailed.Cells.Item(64, "C").Value = "forecast"
$detailed.Cells.Item(64, "D").Value = 46044
$detailed.Cells.Item(64, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(64, "E").Value = "ON"

$detailed.Cells.Item(65, "A").Value = 46044.3125
$detailed.Cells.Item(65, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(65, "B").Value = 0.51
$detailed.Cells.Item(65, "C").Value = "forecast"
$detailed.Cells.Item(65, "D").Value = 46044
$detailed.Cells.Item(65, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(65, "E").Value = "ON"

$detailed.Cells.Item(66, "A").Value = 46044.33333333334
$detailed.Cells.Item(66, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(66, "B").Value = -5.17224
$detailed.Cells.Item(66, "C").Value = "forecast"
$detailed.Cells.Item(66, "D").Value = 46044
$detailed.Cells.Item(66, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(66, "E").Value = "ON"

$detailed.Cells.Item(67, "A").Value = 46044.35416666666
$detailed.Cells.Item(67, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(67, "B").Value = -7.66671
$detailed.Cells.Item(67, "C").Value = "forecast"
$detailed.Cells.Item(67, "D").Value = 46044
$detailed.Cells.Item(67, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(67, "E").Value = "ON"

$detailed.Cells.Item(68, "A").Value = 46044.375
$detailed.Cells.Item(68, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(68, "B").Value = -10
$detailed.Cells.Item(68, "C").Value = "forecast"
$detailed.Cells.Item(68, "D").Value = 46044
$detailed.Cells.Item(68, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(68, "E").Value = "ON"

$detailed.Cells.Item(69, "A").Value = 46044.39583333334
$detailed.Cells.Item(69, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(69, "B").Value = -9.628729999999999
$detailed.Cells.Item(69, "C").Value = "forecast"
$detailed.Cells.Item(69, "D").Value = 46044
$detailed.Cells.Item(69, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(69, "E").Value = "ON"

$detailed.Cells.Item(70, "A").Value = 46044.41666666666
$detailed.Cells.Item(70, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(70, "B").Value = -5.58973
$detailed.Cells.Item(70, "C").Value = "forecast"
$detailed.Cells.Item(70, "D").Value = 46044
$detailed.Cells.Item(70, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(70, "E").Value = "ON"

$detailed.Cells.Item(71, "A").Value = 46044.4375
$detailed.Cells.Item(71, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(71, "B").Value = -7.67552
$detailed.Cells.Item(71, "C").Value = "forecast"
$detailed.Cells.Item(71, "D").Value = 46044
$detailed.Cells.Item(71, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(71, "E").Value = "ON"

$detailed.Cells.Item(72, "A").Value = 46044.45833333334
$detailed.Cells.Item(72, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(72, "B").Value = -7.84586
$detailed.Cells.Item(72, "C").Value = "forecast"
$detailed.Cells.Item(72, "D").Value = 46044
$detailed.Cells.Item(72, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(72, "E").Value = "ON"

$detailed.Cells.Item(73, "A").Value = 46044.47916666666
$detailed.Cells.Item(73, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(73, "B").Value = -5.58973
$detailed.Cells.Item(73, "C").Value = "forecast"
$detailed.Cells.Item(73, "D").Value = 46044
$detailed.Cells.Item(73, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(73, "E").Value = "ON"

$detailed.Cells.Item(74, "A").Value = 46044.5
$detailed.Cells.Item(74, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(74, "B").Value = -7.91485
$detailed.Cells.Item(74, "C").Value = "forecast"
$detailed.Cells.Item(74, "D").Value = 46044
$detailed.Cells.Item(74, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(74, "E").Value = "ON"

$detailed.Cells.Item(75, "A").Value = 46044.52083333334
$detailed.Cells.Item(75, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(75, "B").Value = -7.89793
$detailed.Cells.Item(75, "C").Value = "forecast"
$detailed.Cells.Item(75, "D").Value = 46044
$detailed.Cells.Item(75, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(75, "E").Value = "ON"

$detailed.Cells.Item(76, "A").Value = 46044.54166666666
$detailed.Cells.Item(76, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(76, "B").Value = -6.67846
$detailed.Cells.Item(76, "C").Value = "forecast"
$detailed.Cells.Item(76, "D").Value = 46044
$detailed.Cells.Item(76, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(76, "E").Value = "ON"

$detailed.Cells.Item(77, "A").Value = 46044.5625
$detailed.Cells.Item(77, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(77, "B").Value = -5.50985
$detailed.Cells.Item(77, "C").Value = "forecast"
$detailed.Cells.Item(77, "D").Value = 46044
$detailed.Cells.Item(77, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(77, "E").Value = "ON"

$detailed.Cells.Item(78, "A").Value = 46044.58333333334
$detailed.Cells.Item(78, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(78, "B").Value = -8
$detailed.Cells.Item(78, "C").Value = "forecast"
$detailed.Cells.Item(78, "D").Value = 46044
$detailed.Cells.Item(78, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(78, "E").Value = "ON"

$detailed.Cells.Item(79, "A").Value = 46044.60416666666
$detailed.Cells.Item(79, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(79, "B").Value = -5.50985
$detailed.Cells.Item(79, "C").Value = "forecast"
$detailed.Cells.Item(79, "D").Value = 46044
$detailed.Cells.Item(79, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(79, "E").Value = "ON"

$detailed.Cells.Item(80, "A").Value = 46044.625
$detailed.Cells.Item(80, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(80, "B").Value = -5.50985
$detailed.Cells.Item(80, "C").Value = "forecast"
$detailed.Cells.Item(80, "D").Value = 46044
$detailed.Cells.Item(80, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(80, "E").Value = "ON"

$detailed.Cells.Item(81, "A").Value = 46044.64583333334
$detailed.Cells.Item(81, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(81, "B").Value = 0.51
$detailed.Cells.Item(81, "C").Value = "forecast"
$detailed.Cells.Item(81, "D").Value = 46044
$detailed.Cells.Item(81, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(81, "E").Value = "ON"

$detailed.Cells.Item(82, "A").Value = 46044.66666666666
$detailed.Cells.Item(82, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(82, "B").Value = -1.27042
$detailed.Cells.Item(82, "C").Value = "forecast"
$detailed.Cells.Item(82, "D").Value = 46044
$detailed.Cells.Item(82, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(82, "E").Value = "OFF"

$detailed.Cells.Item(83, "A").Value = 46044.6875
$detailed.Cells.Item(83, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(83, "B").Value = -5.25906
$detailed.Cells.Item(83, "C").Value = "forecast"
$detailed.Cells.Item(83, "D").Value = 46044
$detailed.Cells.Item(83, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(83, "E").Value = "OFF"

$detailed.Cells.Item(84, "A").Value = 46044.70833333334
$detailed.Cells.Item(84, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(84, "B").Value = -12.32215
$detailed.Cells.Item(84, "C").Value = "forecast"
$detailed.Cells.Item(84, "D").Value = 46044
$detailed.Cells.Item(84, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(84, "E").Value = "OFF"

$detailed.Cells.Item(85, "A").Value = 46044.72916666666
$detailed.Cells.Item(85, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(85, "B").Value = -8.505509999999999
$detailed.Cells.Item(85, "C").Value = "forecast"
$detailed.Cells.Item(85, "D").Value = 46044
$detailed.Cells.Item(85, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(85, "E").Value = "OFF"

$detailed.Cells.Item(86, "A").Value = 46044.75
$detailed.Cells.Item(86, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(86, "B").Value = -6
$detailed.Cells.Item(86, "C").Value = "forecast"
$detailed.Cells.Item(86, "D").Value = 46044
$detailed.Cells.Item(86, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(86, "E").Value = "OFF"

$detailed.Cells.Item(87, "A").Value = 46044.77083333334
$detailed.Cells.Item(87, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(87, "B").Value = -2.90964
$detailed.Cells.Item(87, "C").Value = "forecast"
$detailed.Cells.Item(87, "D").Value = 46044
$detailed.Cells.Item(87, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(87, "E").Value = "OFF"

$detailed.Cells.Item(88, "A").Value = 46044.79166666666
$detailed.Cells.Item(88, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(88, "B").Value = 10.11477
$detailed.Cells.Item(88, "C").Value = "forecast"
$detailed.Cells.Item(88, "D").Value = 46044
$detailed.Cells.Item(88, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(88, "E").Value = "OFF"

$detailed.Cells.Item(89, "A").Value = 46044.8125
$detailed.Cells.Item(89, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(89, "B").Value = 57.46185
$detailed.Cells.Item(89, "C").Value = "forecast"
$detailed.Cells.Item(89, "D").Value = 46044
$detailed.Cells.Item(89, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(89, "E").Value = "OFF"

$detailed.Cells.Item(90, "A").Value = 46044.83333333334
$detailed.Cells.Item(90, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(90, "B").Value = 40.46468
$detailed.Cells.Item(90, "C").Value = "forecast"
$detailed.Cells.Item(90, "D").Value = 46044
$detailed.Cells.Item(90, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(90, "E").Value = "OFF"

$detailed.Cells.Item(91, "A").Value = 46044.85416666666
$detailed.Cells.Item(91, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(91, "B").Value = 48.81991
$detailed.Cells.Item(91, "C").Value = "forecast"
$detailed.Cells.Item(91, "D").Value = 46044
$detailed.Cells.Item(91, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(91, "E").Value = "OFF"

$detailed.Cells.Item(92, "A").Value = 46044.875
$detailed.Cells.Item(92, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(92, "B").Value = 51.47526
$detailed.Cells.Item(92, "C").Value = "forecast"
$detailed.Cells.Item(92, "D").Value = 46044
$detailed.Cells.Item(92, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(92, "E").Value = "OFF"

$detailed.Cells.Item(93, "A").Value = 46044.89583333334
$detailed.Cells.Item(93, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(93, "B").Value = 48.40682
$detailed.Cells.Item(93, "C").Value = "forecast"
$detailed.Cells.Item(93, "D").Value = 46044
$detailed.Cells.Item(93, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(93, "E").Value = "OFF"

$detailed.Cells.Item(94, "A").Value = 46044.91666666666
$detailed.Cells.Item(94, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(94, "B").Value = 47.14969
$detailed.Cells.Item(94, "C").Value = "forecast"
$detailed.Cells.Item(94, "D").Value = 46044
$detailed.Cells.Item(94, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(94, "E").Value = "OFF"

$detailed.Cells.Item(95, "A").Value = 46044.9375
$detailed.Cells.Item(95, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(95, "B").Value = 43.519
$detailed.Cells.Item(95, "C").Value = "forecast"
$detailed.Cells.Item(95, "D").Value = 46044
$detailed.Cells.Item(95, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(95, "E").Value = "OFF"

$detailed.Cells.Item(96, "A").Value = 46044.95833333334
$detailed.Cells.Item(96, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(96, "B").Value = 47.17773
$detailed.Cells.Item(96, "C").Value = "forecast"
$detailed.Cells.Item(96, "D").Value = 46044
$detailed.Cells.Item(96, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(96, "E").Value = "OFF"

$detailed.Cells.Item(97, "A").Value = 46044.97916666666
$detailed.Cells.Item(97, "A").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(97, "B").Value = 47.47163
$detailed.Cells.Item(97, "C").Value = "forecast"
$detailed.Cells.Item(97, "D").Value = 46044
$detailed.Cells.Item(97, "D").NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(97, "E").Value = "OFF"

